$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: property_category column I, rows 2-5, change "land" -> "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I5").Value = "building"

# 汽車 (Car) sheet: property_category column H, row 2, change "land" -> "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
